$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-03-19 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-20 Wednesday", 2)

$d.Content.Find.Execute("31×38=1178", $true, $false, $false, $false, $false, $true, 1, $false, "92×39=3588", 2)
$d.Content.Find.Execute("18×34=612", $true, $false, $false, $false, $false, $true, 1, $false, "98×79=7742", 2)
$d.Content.Find.Execute("83×90=7470", $true, $false, $false, $false, $false, $true, 1, $false, "74×60=4440", 2)
$d.Content.Find.Execute("22×14=308", $true, $false, $false, $false, $false, $true, 1, $false, "83×32=2656", 2)
$d.Content.Find.Execute("74×63=4662", $true, $false, $false, $false, $false, $true, 1, $false, "81×53=4293", 2)

$d.Content.Find.Execute("37×35=1295", $true, $false, $false, $false, $false, $true, 1, $false, "20×92=1840", 2)
$d.Content.Find.Execute("73×46=3358", $true, $false, $false, $false, $false, $true, 1, $false, "65×54=3510", 2)
$d.Content.Find.Execute("86×14=1204", $true, $false, $false, $false, $false, $true, 1, $false, "97×35=3395", 2)
$d.Content.Find.Execute("40×76=3040", $true, $false, $false, $false, $false, $true, 1, $false, "60×43=2580", 2)
$d.Content.Find.Execute("36×72=2592", $true, $false, $false, $false, $false, $true, 1, $false, "57×53=3021", 2)

$d.Content.Find.Execute("73×88=6424", $true, $false, $false, $false, $false, $true, 1, $false, "29×15=435", 2)
$d.Content.Find.Execute("84×54=4536", $true, $false, $false, $false, $false, $true, 1, $false, "38×46=1748", 2)
$d.Content.Find.Execute("25×21=525", $true, $false, $false, $false, $false, $true, 1, $false, "91×63=5733", 2)
$d.Content.Find.Execute("36×11=396", $true, $false, $false, $false, $false, $true, 1, $false, "33×61=2013", 2)
$d.Content.Find.Execute("42×36=1512", $true, $false, $false, $false, $false, $true, 1, $false, "38×67=2546", 2)

$d.Content.Find.Execute("48×70=3360", $true, $false, $false, $false, $false, $true, 1, $false, "47×88=4136", 2)
$d.Content.Find.Execute("89×96=8544", $true, $false, $false, $false, $false, $true, 1, $false, "65×22=1430", 2)
$d.Content.Find.Execute("84×77=6468", $true, $false, $false, $false, $false, $true, 1, $false, "55×35=1925", 2)
$d.Content.Find.Execute("21×70=1470", $true, $false, $false, $false, $false, $true, 1, $false, "56×23=1288", 2)
$d.Content.Find.Execute("57×88=5016", $true, $false, $false, $false, $false, $true, 1, $false, "52×60=3120", 2)

$d.Content.Find.Execute("60×42=2520", $true, $false, $false, $false, $false, $true, 1, $false, "43×17=731", 2)
$d.Content.Find.Execute("53×79=4187", $true, $false, $false, $false, $false, $true, 1, $false, "33×51=1683", 2)
$d.Content.Find.Execute("91×64=5824", $true, $false, $false, $false, $false, $true, 1, $false, "53×31=1643", 2)
$d.Content.Find.Execute("17×47=799", $true, $false, $false, $false, $false, $true, 1, $false, "77×38=2926", 2)
$d.Content.Find.Execute("64×67=4288", $true, $false, $false, $false, $false, $true, 1, $false, "18×68=1224", 2)
